$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 2345
    $ws.Range("F3").Value = 1818
    $ws.Range("F4").Value = 351
    $ws.Range("F5").Value = 1121
    $ws.Range("F6").Value = 1027
    $ws.Range("F8").Value = 5912
}
